# Re-order the list of names/emails stored in the "Recorded By" column (G)
# so that entries follow a fixed priority order instead of the order in
# which they were originally appended.
#
# Priority (highest first):
#   admin@admin.com > System > system > dnasr281@gmail.com > backup@backdoor.com
# Any unrecognized name keeps a low priority and is sorted after the known
# ones, preserving its relative position among other unknowns (stable sort).

function Get-RecorderPriority($name) {
    if ($name.Equals("admin@admin.com")) { return 0 }
    if ($name.Equals("System")) { return 1 }
    if ($name.Equals("system")) { return 2 }
    if ($name.Equals("dnasr281@gmail.com")) { return 3 }
    if ($name.Equals("backup@backdoor.com")) { return 4 }
    return 99
}

function Sort-Recorders($value) {
    $parts = @($value -split ", ")
    $n = $parts.Count

    # Stable bubble sort by priority (keeps relative order of equal-priority items)
    for ($i = 0; $i -lt $n; $i++) {
        for ($j = 0; $j -lt ($n - $i - 1); $j++) {
            $p1 = Get-RecorderPriority $parts[$j]
            $p2 = Get-RecorderPriority $parts[$j + 1]
            if ($p1 -gt $p2) {
                $tmp = $parts[$j]
                $parts[$j] = $parts[$j + 1]
                $parts[$j + 1] = $tmp
            }
        }
    }

    return ($parts -join ", ")
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

# Column G = "Recorded By" (column index 7); row 1 is the header row.
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value2
    if ($current -ne $null -and $current -ne "") {
        $updated = Sort-Recorders $current
        if (-not $updated.Equals($current)) {
            $cell.Value = $updated
        }
    }
}
